# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (between the existing "2021-Q4" sheet
# and the "总计" (totals) sheet), populates it with the quarter's fund
# holdings, and prepends a matching summary row to the "总计" sheet.

function Set-TextValue($cell, $val) {
    # Force the cell to be stored as text (even if the string looks like a
    # number, e.g. fund codes or "4.20") while keeping the default style.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetRef = $wb.Worksheets.Item("总计")

# Create the new sheet positioned right before the "总计" sheet so the
# final order is: 2021-Q4, 2022-Q1, 总计.
$newSheet = $wb.Worksheets.Add($totalSheetRef)
$newSheet.Name = "2022-Q1"

# NOTE: worksheet variables in this host track by *position*, not by
# identity, so after inserting a sheet at the "总计" sheet's old slot the
# stale $totalSheetRef variable now actually points at the freshly
# inserted sheet. Re-resolve the "总计" sheet by name once the insert is
# complete so later edits land on the correct sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the cell formatting (borders/fonts/alignment) from the existing
# quarter sheet so the new sheet matches the established layout/styles.
$q4Sheet.Range("A1:H6").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

$rows = @(
    @{ idx = 0; code = "010676"; name = "光大保德信新机遇混合";     scale = "4.20"; pos = "85.57"; pct = "4.29"; mv = "0.1802"; rank = 10 },
    @{ idx = 1; code = "009986"; name = "天弘创新领航混合A";       scale = "2.58"; pos = "73.39"; pct = "2.37"; mv = "0.0611"; rank = 10 },
    @{ idx = 2; code = "002630"; name = "江信瑞福灵活配置混合A";   scale = "0.52"; pos = "43.17"; pct = "2.47"; mv = "0.0128"; rank = 8 },
    @{ idx = 3; code = "002631"; name = "江信瑞福灵活配置混合C";   scale = "0.50"; pos = "43.17"; pct = "2.47"; mv = "0.0124"; rank = 8 },
    @{ idx = 4; code = "009987"; name = "天弘创新领航混合C";       scale = "0.51"; pos = "73.39"; pct = "2.37"; mv = "0.0121"; rank = 10 }
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r,1).Value = $row.idx
    Set-TextValue $newSheet.Cells.Item($r,2) $row.code
    Set-TextValue $newSheet.Cells.Item($r,3) $row.name
    Set-TextValue $newSheet.Cells.Item($r,4) $row.scale
    Set-TextValue $newSheet.Cells.Item($r,5) $row.pos
    Set-TextValue $newSheet.Cells.Item($r,6) $row.pct
    Set-TextValue $newSheet.Cells.Item($r,7) $row.mv
    $newSheet.Cells.Item($r,8).Value = $row.rank
    $r = $r + 1
}

# Insert a new summary row for 2022-Q1 at the top of the "总计" sheet's
# data (row 2), pushing the existing 2021-Q4 row down to row 3.
$totalSheet.Rows.Item(2).Insert()

# Row-insert leaves the new row's cells with an inherited "bold header"
# style; re-apply the bordered index-column style (copied from the row
# that was just pushed down to A3) to A2, and reset B2:D2 back to the
# plain/default style used by ordinary data cells.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Style = "Normal"
Set-TextValue $totalSheet.Cells.Item(2,2) "2022-Q1"
$totalSheet.Cells.Item(2,3).Style = "Normal"
$totalSheet.Cells.Item(2,3).Value = 5
$totalSheet.Cells.Item(2,4).Style = "Normal"
$totalSheet.Cells.Item(2,4).Value = 0.28

# Renumber the pre-existing 2021-Q4 row's index cell (A3) to 1, matching
# the diff (rows are 0-indexed sequentially down the sheet).
$totalSheet.Cells.Item(3,1).Value = 1

# Keep the originally active sheet/tab selected (unchanged by this edit).
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q4Sheet.Select()
